$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.953.76"
$ws.Range("E2").Value = "  -2.19%  "

# Row 3
$ws.Range("D3").Value = "3.746.58"
$ws.Range("E3").Value = "  -0.79%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.78%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.43%  "

# Row 7
$ws.Range("D7").Value = "3.745.63"
$ws.Range("E7").Value = "  -0.78%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  -2.00%  "

# Row 10
$ws.Range("E10").Value = "  -3.63%  "

# Row 11
$ws.Range("E11").Value = "  -2.64%  "

# Row 12
$ws.Range("E12").Value = "  -0.87%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.26%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.39%  "

# Row 15
$ws.Range("D15").Value = "4.375.71"
$ws.Range("E15").Value = "  -0.85%  "

# Row 16
$ws.Range("D16").Value = "3.746.79"
$ws.Range("E16").Value = "  -0.47%  "

# Row 17
$ws.Range("D17").Value = "66.992.05"
$ws.Range("E17").Value = "  -2.18%  "

# Row 18
$ws.Range("E18").Value = "  -2.85%  "

# Row 19
$ws.Range("E19").Value = "  -0.09%  "

# Row 20
$ws.Range("E20").Value = "  -2.23%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.70%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "452.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.33%  "

# Row 23
$ws.Range("E23").Value = "  -1.94%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000146"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.54%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.92%  "

# Row 26
$ws.Range("E26").Value = "  -5.47%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.02%  "

# Row 28
$ws.Range("E28").Value = "  +0.05%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.73%  "

# Row 30
$ws.Range("E30").Value = "  -1.98%  "

# Row 31
$ws.Range("E31").Value = "  -3.41%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.23%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.32%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.61%  "

# Row 35
$ws.Range("E35").Value = "  +0.08%  "

# Row 36
$ws.Range("D36").Value = "3.700.01"
$ws.Range("E36").Value = "  -0.85%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0987"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.14%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.137"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.98%  "

# Row 39
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.36%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.986"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.77%  "

# Row 41
$ws.Range("E41").Value = "  -2.45%  "

# Row 42
$ws.Range("E42").Value = "  -0.07%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.21%  "

# Row 45
$ws.Range("E45").Value = "  -4.66%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.63"
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = "  -3.26%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "146.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.71%  "

# Row 49
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "385.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.05%  "

# Row 50
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.87%  "

# Row 51
$ws.Range("D51").Value = "2.735.93"
$ws.Range("E51").Value = "  +1.66%  "
